$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.167.94'
$ws.Range("E2").Value = '  -2.98%  '

$ws.Range("D3").Value = '1.607.17'
$ws.Range("E3").Value = '  -2.45%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").Value = '302.15'
$ws.Range("E6").Value = '  -2.12%  '

$ws.Range("D7").Value = '0.3764'
$ws.Range("E7").Value = '  -3.19%  '

$ws.Range("D8").Value = '0.3637'
$ws.Range("E8").Value = '  -4.75%  '

$ws.Range("D9").Value = '48.63'
$ws.Range("E9").Value = '  -5.14%  '

$ws.Range("D10").Value = '1.003'
$ws.Range("E10").Value = '  -0.06%  '

$ws.Range("D11").Value = '1.264'
$ws.Range("E11").Value = '  -6.03%  '

$ws.Range("D12").Value = '0.08052'
$ws.Range("E12").Value = '  -4.33%  '

$ws.Range("D13").Value = '22.89'
$ws.Range("E13").Value = '  -3.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.560'
$ws.Range("E14").Value = '  -7.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.670'
$ws.Range("E15").Value = '  -2.20%  '

$ws.Range("D16").Value = '0.00001261'
$ws.Range("E16").Value = '  -3.83%  '

$ws.Range("D17").Value = '1.602.93'
$ws.Range("E17").Value = '  -2.96%  '

$ws.Range("D18").Value = '91.38'
$ws.Range("E18").Value = '  -3.07%  '

$ws.Range("D19").Value = '0.06788'
$ws.Range("E19").Value = '  -2.84%  '

$ws.Range("D20").Value = '18.26'
$ws.Range("E20").Value = '  -7.09%  '

$ws.Range("D21").Value = '6.536'
$ws.Range("E21").Value = '  -5.44%  '

$ws.Range("D22").Value = '1.004'
$ws.Range("E22").Value = '  +0.15%  '

$ws.Range("D23").Value = '13.04'
$ws.Range("E23").Value = '  -4.63%  '

$ws.Range("D24").Value = '23.191.12'
$ws.Range("E24").Value = '  -2.92%  '

$ws.Range("D25").Value = '2.353'
$ws.Range("E25").Value = '  -4.16%  '

$ws.Range("D26").Value = '2.878'
$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("D27").Value = '21.02'
$ws.Range("E27").Value = '  -4.49%  '

$ws.Range("D28").Value = '150.28'
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("D29").Value = '5.266'
$ws.Range("E29").Value = '  -2.80%  '

$ws.Range("D30").Value = '132.16'
$ws.Range("E30").Value = '  -4.46%  '

$ws.Range("D31").Value = '2.385'
$ws.Range("E31").Value = '  -4.53%  '

$ws.Range("D32").Value = '6.787'
$ws.Range("E32").Value = '  -12.74%  '

$ws.Range("D33").Value = '1.780.78'
$ws.Range("E33").Value = '  -2.83%  '

$ws.Range("D34").Value = '0.9639'
$ws.Range("E34").Value = '  -7.61%  '

$ws.Range("D35").Value = '0.07708'
$ws.Range("E35").Value = '  -4.26%  '

$ws.Range("D36").Value = '0.02768'
$ws.Range("E36").Value = '  -6.21%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '6.225'
$ws.Range("E37").Value = '  -7.50%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = '0.2532'
$ws.Range("E38").Value = '  -5.27%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.10'
$ws.Range("E39").Value = '  -6.70%  '

$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.08839'
$ws.Range("E40").Value = '  -2.93%  '

$ws.Range("D41").Value = '1.392'
$ws.Range("E41").Value = '  -1.88%  '

$ws.Range("D42").Value = '0.7132'
$ws.Range("E42").Value = '  -5.30%  '

$ws.Range("D43").Value = '12.77'
$ws.Range("E43").Value = '  -4.84%  '

$ws.Range("D44").Value = '15.65'
$ws.Range("E44").Value = '  -3.81%  '

$ws.Range("D45").Value = '0.6568'
$ws.Range("E45").Value = '  -5.20%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").Value = '2.288'
$ws.Range("E47").Value = '  -6.62%  '

$ws.Range("D48").Value = '3.976'
$ws.Range("E48").Value = '  -2.60%  '

$ws.Range("D49").Value = '0.07983'

$ws.Range("D50").Value = '131.77'
$ws.Range("E50").Value = '  -1.62%  '

$ws.Range("D51").Value = '1.167'
$ws.Range("E51").Value = '  -3.23%  '

